$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price column so numeric-looking strings
# (e.g. "615.26", "1.00", thousand-dot formatted numbers) are not
# auto-converted to floating point numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "64.730.07"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "3.155.39"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "615.26"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").Value = "147.98"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.152.05"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("D11").Value = "5.47"
$ws.Range("E11").Value = "  -2.82%  "
$ws.Range("D12").Value = "0.473"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "35.82"
$ws.Range("E14").Value = "  -3.44%  "
$ws.Range("D15").Value = "3.674.05"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("D17").Value = "64.701.10"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "3.156.83"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "6.94"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").Value = "481.78"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").Value = "14.73"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "8.04"
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("D24").Value = "13.79"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").Value = "84.43"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -3.22%  "
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("D29").Value = "7.02"
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").Value = "0.118"
$ws.Range("E30").Value = "  -2.36%  "
$ws.Range("D31").Value = "2.09"
$ws.Range("E31").Value = "  -7.74%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "2.72"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").Value = "26.51"
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("D36").Value = "0.0₃0783"
$ws.Range("E36").Value = "  +4.05%  "
$ws.Range("D37").Value = "6.02"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("D39").Value = "52.95"
$ws.Range("E39").Value = "  -3.25%  "
$ws.Range("D40").Value = "460.75"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("E42").Value = "  -5.26%  "
$ws.Range("D43").Value = "8.41"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("D44").Value = "2.845.55"
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("E45").Value = "  -4.40%  "
$ws.Range("D46").Value = "0.269"
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("D47").Value = "2.48"
$ws.Range("E47").Value = "  +5.02%  "
$ws.Range("D48").Value = "26.63"
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").Value = "120.75"
$ws.Range("E51").Value = "  +0.48%  "

# Restore the original (default) style so no stray number format/style
# index is left behind on these cells.
$priceRange.Style = "Normal"
